$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.147.23"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -4.34%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.652.94"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -3.42%  "

# Row 4
$ws.Range("E4").Value = "  +0.09%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.86"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -3.73%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5114"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -3.40%  "

# Row 7
$ws.Range("E7").Value = "  +0.05%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2584"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -3.05%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06430"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -3.52%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.98"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -4.12%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07802"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +1.41%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.658.69"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -3.25%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.283"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -4.75%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.880.43"
$ws.Range("D14").ClearFormats()

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5521"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -5.24%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅8037"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -2.09%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.11"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -5.74%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.160.24"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -4.37%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.006"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.05%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "210.81"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -4.66%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.411"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -4.66%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.08"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -3.44%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.034"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.52%  "

# Row 24
$ws.Range("E24").Value = "  +0.01%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.77"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.68%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.738"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +2.88%  "

# Row 27
$ws.Range("E27").Value = "  -1.80%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.984"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -3.38%  "

# Row 29
$ws.Range("E29").Value = "  -2.57%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05114"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -4.10%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.246"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -3.40%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.343"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -3.43%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.220"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -6.19%  "

# Row 34
$ws.Range("E34").Value = "  -4.16%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.747"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -4.56%  "

# Row 36
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9264"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -2.58%  "

# Row 37
$ws.Range("B37").Value = "HuobiToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.360"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -1.43%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.164.91"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +4.57%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5687"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -2.49%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01585"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -2.88%  "

# Row 41
$ws.Range("B41").Value = "mCoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.553"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.57%  "

# Row 42
$ws.Range("B42").Value = "PaxDollar"
$ws.Range("C42").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.005"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.03%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8297"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.05%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.652"
$ws.Range("D44").ClearFormats()

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "100.37"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.82%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.790.79"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -3.41%  "

# Row 47
$ws.Range("E47").Value = "  -0.08%  "

# Row 48
$ws.Range("E48").Value = "  +0.12%  "

# Row 49
$ws.Range("E49").Value = "  -3.45%  "

# Row 50
$ws.Range("E50").Value = "  +0.11%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.862"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -2.73%  "
